# Regenerate the handback report timestamps (commit: "Generate Report for Handback").
# Several "...Xliff Generate Date" / "...Handback DateTime" columns get refreshed
# with a later run time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets("Overview")
$wsZhCn     = $wb.Worksheets("zh-cn")
$wsDeDe     = $wb.Worksheets("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-16 06:58:36"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file.
$wsZhCn.Range("H2").Value = "2016-08-16 06:58:31"
$wsZhCn.Range("K2").Value = "2016-08-16 06:58:50"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file.
$wsDeDe.Range("H2").Value = "2016-08-16 06:58:36"
$wsDeDe.Range("K2").Value = "2016-08-16 06:58:57"
